# Typo corrections in fellows ppt
# 1) Slide 12, "Rectangle 8" shape: "pretect" -> "protect"
# 2) Slide 9, "TextBox 2" shape: "visulazation" -> "visualization"

$p = $ppt.ActivePresentation

# --- Fix 1: slide 12 ---
$s12 = $p.Slides.Item(12)
$shape12 = $s12.Shapes.Item(7)
$tr12 = $shape12.TextFrame.TextRange
$para6 = $tr12.Paragraphs(6, 1)
$bad = $para6.Characters(32, 8)   # "pretect " (includes trailing space)
$bad.Text = "protect "

# --- Fix 2: slide 9 ---
$s9 = $p.Slides.Item(9)
$shape9 = $s9.Shapes.Item(2)
$tr9 = $shape9.TextFrame.TextRange
$para9 = $tr9.Paragraphs(9, 1)
$bad2 = $para9.Characters(39, 13)  # "visulazation " (includes trailing space)
$bad2.Text = "visualization "
